$rowData = @{
    2 = @{ "G"=54.53585066666667; "H"=163.607552; "I"=0.3031388658437607; "J"=0.3031388658437607; "M"=9.423852333333334; "N"=28.271557; "O"=0.06654336290212845; "P"=0.06654336290212845; "Q"=513.937803555385; "R"=4625.440231998464; "S"=0.020171879559581; "T"=0.020171879559581 }
    3 = @{ "G"=54.53585066666667; "H"=163.607552; "I"=0.3031388658437607; "J"=0.3031388658437607; "M"=50.59256466666667; "N"=151.777694; "O"=0.3572423751649123; "P"=0.3572423751649123; "Q"=2759.10855150501; "R"=24831.97696354509; "S"=0.1082940484388228; "T"=0.1082940484388228 }
    4 = @{ "G"=54.53585066666667; "H"=163.607552; "I"=0.3031388658437607; "J"=0.3031388658437607; "M"=26.84076266666667; "N"=80.522288; "O"=0.1895270158659356; "P"=0.1895270158659356; "Q"=1463.783824568775; "R"=13174.05442111898; "S"=0.05745300463635214; "T"=0.05745300463635216 }
    5 = @{ "G"=54.53585066666667; "H"=163.607552; "I"=0.3031388658437607; "J"=0.3031388658437607; "M"=54.762539; "N"=164.287617; "O"=0.3866872460670236; "P"=0.3866872460670236; "Q"=2986.521649031509; "R"=26878.69484128358; "S"=0.1172199332090047; "T"=0.1172199332090048 }
    6 = @{ "G"=18.46467533333333; "H"=55.394026; "I"=0.1026363515063155; "J"=0.1026363515063155; "M"=9.423852333333334; "N"=28.271557; "O"=0.06654336290212845; "P"=0.06654336290212845; "Q"=174.0083737242758; "R"=1566.075363518482; "S"=0.00682976798523517; "T"=0.006829767985235171 }
    7 = @{ "G"=18.46467533333333; "H"=55.394026; "I"=0.1026363515063155; "J"=0.1026363515063155; "M"=50.59256466666667; "N"=151.777694; "O"=0.3572423751649123; "P"=0.3572423751649123; "Q"=934.1752808506715; "R"=8407.577527656043; "S"=0.03666605399037697; "T"=0.03666605399037698 }
    8 = @{ "G"=18.46467533333333; "H"=55.394026; "I"=0.1026363515063155; "J"=0.1026363515063155; "M"=26.84076266666667; "N"=80.522288; "O"=0.1895270158659356; "P"=0.1895270158659356; "Q"=495.6059683390542; "R"=4460.453715051488; "S"=0.0194523614203592; "T"=0.0194523614203592 }
    9 = @{ "G"=18.46467533333333; "H"=55.394026; "I"=0.1026363515063155; "J"=0.1026363515063155; "M"=54.762539; "N"=164.287617; "O"=0.3866872460670236; "P"=0.3866872460670236; "Q"=1011.172503064005; "R"=9100.552527576041; "S"=0.03968816811034415; "T"=0.03968816811034415 }
    10 = @{ "G"=12.55635966666667; "H"=37.669079; "I"=0.06979483370938171; "J"=0.06979483370938172; "M"=9.423852333333334; "N"=28.271557; "O"=0.06654336290212845; "P"=0.06654336290212845; "Q"=118.3292793428892; "R"=1064.963514086003; "S"=0.004644382948217095; "T"=0.004644382948217096 }
    11 = @{ "G"=12.55635966666667; "H"=37.669079; "I"=0.06979483370938171; "J"=0.06979483370938172; "M"=50.59256466666667; "N"=151.777694; "O"=0.3572423751649123; "P"=0.3572423751649123; "Q"=635.2584384137584; "R"=5717.325945723825; "S"=0.02493367216857961; "T"=0.02493367216857961 }
    12 = @{ "G"=12.55635966666667; "H"=37.669079; "I"=0.06979483370938171; "J"=0.06979483370938172; "M"=26.84076266666667; "N"=80.522288; "O"=0.1895270158659356; "P"=0.1895270158659356; "Q"=337.0222697703057; "R"=3033.200427932752; "S"=0.01322800655579832; "T"=0.01322800655579833 }
    13 = @{ "G"=12.55635966666667; "H"=37.669079; "I"=0.06979483370938171; "J"=0.06979483370938172; "M"=54.762539; "N"=164.287617; "O"=0.3866872460670236; "P"=0.3866872460670236; "Q"=687.6181359438602; "R"=6188.563223494742; "S"=0.02698877203678668; "T"=0.02698877203678669 }
    14 = @{ "G"=94.34696966666667; "H"=283.040909; "I"=0.524429948940542; "J"=0.5244299489405421; "M"=9.423852333333334; "N"=28.271557; "O"=0.06654336290212845; "P"=0.06654336290212845; "Q"=889.1119102361459; "R"=8002.007192125313; "S"=0.03489733240909518; "T"=0.03489733240909518 }
    15 = @{ "G"=94.34696966666667; "H"=283.040909; "I"=0.524429948940542; "J"=0.5244299489405421; "M"=50.59256466666667; "N"=151.777694; "O"=0.3572423751649123; "P"=0.3572423751649123; "Q"=4773.255163964872; "R"=42959.29647568385; "S"=0.1873486005671329; "T"=0.1873486005671329 }
    16 = @{ "G"=94.34696966666667; "H"=283.040909; "I"=0.524429948940542; "J"=0.5244299489405421; "M"=26.84076266666667; "N"=80.522288; "O"=0.1895270158659356; "P"=0.1895270158659356; "Q"=2532.344621142199; "R"=22791.10159027979; "S"=0.09939364325342588; "T"=0.09939364325342591 }
    17 = @{ "G"=94.34696966666667; "H"=283.040909; "I"=0.524429948940542; "J"=0.5244299489405421; "M"=54.762539; "N"=164.287617; "O"=0.3866872460670236; "P"=0.3866872460670236; "Q"=5166.67960590265; "R"=46500.11645312385; "S"=0.202790372710888; "T"=0.202790372710888 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}
